# Apply the "EditarProductoData" fix: refresh the base URL, fix the
# discount/price values, rewrite the second test-case row to be a
# validation scenario instead of a duplicate of the first, and drop the
# now-unused "NuevoPrecio" column (K) -- its old neighbour "Valor
# Esperado" (column L) slides left to take its place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -- existing case, base URL + discount + price values refreshed.
$ws.Range("B2").Value = "http://localhost:3000"
$ws.Range("G2").Value = 20
$ws.Range("J2").Value = "Si"

# Row 3 -- now a distinct validation test case.
$ws.Range("B3").Value = "http://localhost:3000"
$ws.Range("G3").Value = " "
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = "Extranjera"
$ws.Range("J3").Value = "No"
$ws.Range("L3").Value = "Rellenar el campo precio correctamente"

# Match the border style used by the rest of the data rows for the Url
# column (it previously had its own slightly different "applyFill" style).
$ws.Range("B2").Borders.LineStyle = 1
$ws.Range("B3").Borders.LineStyle = 1

# Remove the obsolete "NuevoPrecio" column -- "Valor Esperado" shifts
# left from L into K.
$ws.Range("K1:K3").EntireColumn.Delete()

$ws.Range("K2").Select()
